# Generate Report for Handback
# Row 2 in each language sheet (zh-cn, de-de) corresponds to file
# 7e55884b-a229-45ff-8e3b-859d31913147.md, which has now been handed back
# and is in sync with en-US. Update status, populate the "Latest Target
# File" / "Latest Handback File" columns (F/G) with the handback file
# info (mirroring columns A/D, hyperlinks included), and stamp the
# handback datetime.

$wb = $excel.ActiveWorkbook

$langSheets = @(
    @{
        Name = "zh-cn"
        HandbackDateTime = "2016-03-24 20:46:54"
        SourceUrl = "https://github.com/OpenLocalizationTest/oltest/blob/0852d0d04a9811c6d5b749f5238065ad20ab928d/e2e/7e55884b-a229-45ff-8e3b-859d31913147.md"
        HandoffUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bb956c40de845595b930e47f86fd2103817c1566/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/7e55884b-a229-45ff-8e3b-859d31913147.1866ce92fe4c8446a9c7ca037fae0b33c95adb86.zh-cn.xlf"
    },
    @{
        Name = "de-de"
        HandbackDateTime = "2016-03-24 20:47:02"
        SourceUrl = "https://github.com/OpenLocalizationTest/oltest/blob/0852d0d04a9811c6d5b749f5238065ad20ab928d/e2e/7e55884b-a229-45ff-8e3b-859d31913147.md"
        HandoffUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/64523d157b2a3848352a3d2e21813670ad9f0ead/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/7e55884b-a229-45ff-8e3b-859d31913147.1866ce92fe4c8446a9c7ca037fae0b33c95adb86.de-de.xlf"
    }
)

foreach ($info in $langSheets) {
    $ws = $wb.Worksheets.Item($info.Name)

    # Status: "Ready for handoff" -> "Handed back: in sync with en-US"
    $ws.Range("C2").Value = "Handed back: in sync with en-US"

    # Latest Target File (F2) mirrors Source File Name (A2), including
    # its hyperlink to the source .md file.
    $sourceDisplay = $ws.Range("A2").Value()
    $ws.Hyperlinks.Add($ws.Range("F2"), $info.SourceUrl, "", "", $sourceDisplay)

    # Latest Handback File (G2) mirrors Latest Handoff File (D2),
    # including its hyperlink to the handoff .xlf file.
    $handoffDisplay = $ws.Range("D2").Value()
    $ws.Hyperlinks.Add($ws.Range("G2"), $info.HandoffUrl, "", "", $handoffDisplay)

    # Latest Handback DateTime (H2)
    $ws.Range("H2").Value = $info.HandbackDateTime
}
